$d = $word.ActiveDocument

# The document currently ends with a single empty ListParagraph-styled
# paragraph that carries the "_GoBack" bookmark. We need to:
#   1. Insert two brand-new list paragraphs *before* it (questions 1 & 2).
#   2. Add the text for question 3 into that original paragraph, before
#      the bookmark, which must remain at the very end of the document.
#
# NOTE: a Range object captured once and reused across several edits can
# silently grow as text/paragraphs are inserted at its boundaries, so we
# always re-fetch a fresh Range/Paragraph right before each edit below.

# --- New paragraph 1: "After the cells were lysed..." ---
$anchor1 = $d.Paragraphs.Last.Range
$anchor1.InsertParagraphBefore()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p1.Range.InsertBefore("After the cells were lysed, an extensive amount of centrifuging and experiments were performed on the cell-free state. Is there any concern with stability and/or functionality of the components after multiple days of handling? Was anything done to account for this or make sure there was no function/structure loss in the components?")

# --- New paragraph 2: "It seems that there are no details on trial numbers..." ---
$anchor2 = $d.Paragraphs.Last.Range
$anchor2.InsertParagraphBefore()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$r2 = $p2.Range
$r2.InsertAfter("It seems that there are no details on trial numbers. How was statistical significance obtained or rep")
$r2.InsertAfter("orted in the results")
$r2.InsertAfter(" (e.g.")
$r2.InsertAfter(" the difference in protein production")
$r2.InsertAfter(" with or without ")
$r2.InsertAfter("polyuridylic")
$r2.InsertAfter(" acid")
$r2.InsertAfter(" in Table 6")
$r2.InsertAfter(")?")

# --- Question 3 text goes into the original last paragraph, before the bookmark ---
$finalRange = $d.Paragraphs.Last.Range
$finalRange.InsertBefore("Why does increasing soluble RNA")
$finalRange.InsertAfter(" cause the protein production to increase at all (assuming this is significant data)?")
